$d = $word.ActiveDocument

# --- Simple whole-paragraph XML replacements (paragraph count unchanged) ---

# In addition to GO analysis paragraph (Encyclopedia / signaling spellcheck marks)
$frag20 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">In addition to GO analysis, Kyoto </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Encyclopedia</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> of Genes and Genomes (KEGG) pathway analysis was performed to identify </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>signaling</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> pathways that were significantly enriched among the DEGs.</w:t></w:r></w:p>'
$d.Paragraphs(20).Range.InsertXML($frag20)

# Functional enrichment analysis body paragraph (clusterProfiler spellcheck mark)
$frag18 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">To better understand the biological relevance of the differentially expressed genes, functional enrichment analysis was conducted using the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>clusterProfiler</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> package in R. This analysis focused on identifying over</w:t></w:r><w:r><w:t>represented biological processes, molecular functions, and cellular components through Gene Ontology (GO) analysis</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>'
$d.Paragraphs(18).Range.InsertXML($frag18)

# Functional enrichment analysis heading (remove lastRenderedPageBreak)
$frag17 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/></w:rPr><w:t>Functional enrichment analysis</w:t></w:r></w:p>'
$d.Paragraphs(17).Range.InsertXML($frag17)

# Visualization body paragraph (pheatmap / tumor spellcheck marks)
$frag15 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Several visualization techniques were employed to facilitate the interpretation of the results. Volcano plots were generated using </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>ggplot2</w:t></w:r><w:r><w:t xml:space="preserve"> to provide an overview of significantly upregulated and downregulated genes, highlighting their statistical significance. </w:t></w:r><w:r><w:t xml:space="preserve">Heatmaps, created using the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pheatmap</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> package, were used to visualize hierarchical clustering </w:t></w:r><w:r><w:t>patterns</w:t></w:r><w:r><w:t xml:space="preserve"> of DEGs across </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tumor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and normal samples, helping to identify potential gene expression </w:t></w:r><w:r><w:t>identification</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">MA plots were </w:t></w:r><w:r><w:t>constructed</w:t></w:r><w:r><w:t xml:space="preserve"> to </w:t></w:r><w:r><w:t xml:space="preserve">display log2 fold-change values against mean expression levels, assisting in the identification of expression trends. </w:t></w:r></w:p>'
$d.Paragraphs(15).Range.InsertXML($frag15)

# Ellipsis paragraph after Differential expression analysis (add lastRenderedPageBreak)
$frag12 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t>…</w:t></w:r></w:p>'
$d.Paragraphs(12).Range.InsertXML($frag12)

# Differential expression analysis body paragraph (tumor / DESeq spellcheck marks)
$frag11 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">To identify genes that were significantly differentially expressed between </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tumor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and normal lung tissues, differential expression analysis was performed using the DESeq2 package in R. the first step in this analysis involved converting the raw count </w:t></w:r><w:r><w:t>matrix</w:t></w:r><w:r><w:t xml:space="preserve"> into a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>DESeq</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> dataset</w:t></w:r><w:r><w:t xml:space="preserve">, incorporating sample metadata to specify experimental conditions. </w:t></w:r></w:p>'
$d.Paragraphs(11).Range.InsertXML($frag11)

# Computational environment paragraph (merge run + pheatmap/clusterProfiler spellcheck marks)
$frag3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>To ensure reproducibility and streamline the computational workflow, a dedicated environment was created with all group mem</w:t></w:r><w:r><w:t xml:space="preserve">bers. This environment included dependencies such as R </w:t></w:r><w:r><w:t>(</w:t></w:r><w:r><w:t>version 4.2.2</w:t></w:r><w:r><w:t>)</w:t></w:r><w:r><w:t xml:space="preserve">. Whitin R, several packages were installed to facilitate data processing and analysis, including </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>DESeq2</w:t></w:r><w:r><w:t xml:space="preserve"> for differential expression analysis, </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>ggplot2</w:t></w:r><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>pheatmap</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for visualization</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>clusterProfiler</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> for functional enrichment analysis. T</w:t></w:r><w:r><w:t>he entire project was managed through GitHub, ensuring systematic version control</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> allowing multiple contributors to </w:t></w:r><w:r><w:t>work collaboratively while maintaining an organized storage</w:t></w:r><w:r><w:t xml:space="preserve"> of scripts and results. </w:t></w:r></w:p>'
$d.Paragraphs(3).Range.InsertXML($frag3)

# Data acquisition intro paragraph (tumor / GEOquery spellcheck marks)
$frag2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">To investigate differentially expressed genes (DEGs) in lung cancer, RNA sequencing data was obtained from the Gene Expression Omnibus (GEO) database, specifically dataset GSE81089. This dataset consists of transcriptomic profiles of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tumor</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and adjacent normal lung tissues from patients with non-small cell lung cancer (NSCLC). The dataset was accessed and retrieved using the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GEOquery</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> package in R, which allows for the direct downloading and processing of GEO</w:t></w:r><w:r><w:t xml:space="preserve"> datasets. </w:t></w:r></w:p>'
$d.Paragraphs(2).Range.InsertXML($frag2)

# --- Split paragraph 7 and insert the new normalisation discussion paragraphs ---
$p7 = $d.Paragraphs(7)
$p7r = $p7.Range
$searchRange = $d.Range($p7r.Start, $p7r.End)
$null = $searchRange.Find.Execute("structure the dataset. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $searchRange.End
$restEnd = $p7r.End - 1
$d.Range($splitPoint, $restEnd).Delete()

$p7b = $d.Paragraphs(7)
$insertPoint = $d.Range($p7b.Range.End - 1, $p7b.Range.End - 1)
$frag7insert = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t xml:space="preserve">Exploratory analyses were performed </w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t xml:space="preserve">using count data corrected by Fragments Per Kilobase </w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t xml:space="preserve">of transcript per </w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t xml:space="preserve">Million </w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t xml:space="preserve">reads </w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t>mapped (FPKM). Normalisation methods are used to account for technical variabilities such</w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t xml:space="preserve"> as </w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t xml:space="preserve">sequencing </w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t>depth, transcript</w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t xml:space="preserve"> length</w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t>sample-to-sample variability</w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t xml:space="preserve"> and batch-to-batch variability</w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t xml:space="preserve"> (Conesa et. Al, 2016)</w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t xml:space="preserve"> FPKM corrects for variations in</w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t xml:space="preserve">both </w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t>gene length and sequencing depth. It is calculated by … Because it normalises reads by the total number of fragments mapped, it is indicated to compare gene expression levels between different samples. This is not the case with TPM which is useful when comparing gene expression levels within a sample.  -&gt; RPKM?</w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t>..</w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t xml:space="preserve">The statistical tool DESeq2 automatically uses the median-of-ratios normalisation method thus, to avoid distortion of the normalisation process, raw counts data </w:t></w:r><w:r><w:rPr><w:color w:val="7F7F7F" w:themeColor="text1" w:themeTint="80"/></w:rPr><w:t>was provided as input.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The raw count matrix and associated metadata were loaded into R, ensuring that sample identifiers were correctly matched between the two files. Genes with very low expression levels, defined as having an average count below 10 raw reads, were removed to reduce noise and improve statistical power. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$insertPoint.InsertXML($frag7insert)

Write-Host "Edit complete."
